# Update FantaSPL_Classifica sheet with "First Game results - official"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: score update only
$ws.Range("C3").Value = 73

# Row 5
$ws.Range("B5").Value = "Omanta"
$ws.Range("C5").Value = 61

# Row 6
$ws.Range("B6").Value = "CHIAVO VERONA"
$ws.Range("C6").Value = 61

# Row 7
$ws.Range("B7").Value = "Beverly INPS"
$ws.Range("C7").Value = 61

# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Rahal Madrid"
$ws.Range("C8").Value = 59

# Row 9
$ws.Range("B9").Value = "T'eamCulo"
$ws.Range("C9").Value = 55

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "AC Tua"
$ws.Range("C10").Value = 55

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Fel Lazio"
$ws.Range("C11").Value = 55

# Row 12
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Football Meta Academy"

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Non è la seconda squadra di Mazzu, è la prima"
$ws.Range("C13").Value = 54

# Row 14
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "SPL Solo Per Letette"
$ws.Range("C14").Value = 54

# Row 15
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Mazzu è ok"
$ws.Range("C15").Value = 54

# Row 16
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Rapid Viennetta"
$ws.Range("C16").Value = 53

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Calabria Saudita"
$ws.Range("C17").Value = 53

# Row 18
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Limonta United"
$ws.Range("C18").Value = 52

# Row 19
$ws.Range("B19").Value = "Artificially Degenerated"
$ws.Range("C19").Value = 51

# Row 20
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Sesso Paperoga Lamborghini"
$ws.Range("C20").Value = 49

# Row 21
$ws.Range("B21").Value = "Affori Grizzlies"
$ws.Range("C21").Value = 44

# Row 22
$ws.Range("B22").Value = "Si è girato Mazzoud"
$ws.Range("C22").Value = 43

# Row 23
$ws.Range("B23").Value = "Slayer FC"
$ws.Range("C23").Value = 42

# Row 24
$ws.Range("A24").Value = 23

# Row 25
$ws.Range("B25").Value = "Aldo Ritmo"
$ws.Range("C25").Value = 39

# Row 26
$ws.Range("B26").Value = "BaffoImpregnato"
$ws.Range("C26").Value = 38

# Row 27
$ws.Range("B27").Value = "NonCiCapiscoNaMazza"
$ws.Range("C27").Value = 37
